# Swap the data of row 2 and row 3 (keeping columns that are identical
# between the two rows untouched). This corresponds to the two sighting
# records (Spillkråka / Tjäder) trading places in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values differ between row 2 and row 3 and must be swapped.
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R", "Z", "AB")

foreach ($col in $cols) {
    $addr2 = "$col`2"
    $addr3 = "$col`3"
    $v2 = $ws.Range($addr2).Value2
    $v3 = $ws.Range($addr3).Value2
    $ws.Range($addr2).Value = $v3
    $ws.Range($addr3).Value = $v2
}

# AC column: value moves from row 3 to row 2 (row 3's AC becomes empty).
$ws.Range("AC2").Value = "2 tuppar"
$ws.Range("AC3").Value = ""
